# Commit 18: Auth frontend, requests and controllers
#
# This script updates the "Development Phases" worksheet:
#  - D51 keeps showing "Statistic filter" (its shared-string slot just
#    gets reshuffled upstream, no visible change),
#  - Row 52 (D/E/F) is filled in: "Create, Edit for Auth" / DONE / Frontend,
#  - 8 new rows (53-60) are appended describing the auth frontend/backend
#    work,
#  - The used range grows from C2:F52 to C2:F60,
#  - The view no longer freezes on topLeftCell A22 and the new selection
#    sits on D60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Format the rows first (C/E/F numbers & reused text are all already
#    present in the shared-string table, so formatting passes do not
#    influence new shared-string ordering).
# ---------------------------------------------------------------------

# Row 52 sheds the special "last row" formatting it used to carry now
# that the table keeps going - re-use row 51's plain formatting.
$ws.Range("D51:F51").Copy() | Out-Null
$ws.Range("D52:F52").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("E52").Value2 = "DONE"

$rowInfo = @{
    53 = @{ C = 31.1; E = "DONE"; F = "Frontend" }
    54 = @{ C = 31.2; E = "DONE"; F = "Frontend" }
    55 = @{ C = 32;   E = "DONE"; F = "Frontend" }
    56 = @{ C = 33;   E = "DONE"; F = "Frontend" }
    57 = @{ C = 34;   E = "DONE"; F = "Backend"  }
    58 = @{ C = 35;   E = $null;  F = "Backend"  }
    59 = @{ C = 36;   E = $null;  F = "Frontend" }
    60 = @{ C = 37;   E = $null;  F = "Frontend" }
}

foreach ($rowNum in 53..60) {
    $info = $rowInfo[$rowNum]

    # Duplicate the normal data-row formatting (border, font, no fill)
    # from row 50 into every column of the new row.
    $ws.Range("C50:F50").Copy() | Out-Null
    $ws.Range("C" + $rowNum + ":F" + $rowNum).PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false

    $ws.Range("C$rowNum").Value2 = $info.C
    if ($info.E) {
        $ws.Range("E$rowNum").Value2 = $info.E
    }
    $ws.Range("F$rowNum").Value2 = $info.F
}

# Row 60 is the new final row of the table. Give C60 the bold look that
# used to mark the final row of the table.
$ws.Range("C60").Font.Bold = $true

# ---------------------------------------------------------------------
# 2. Now fill in the D column text in the exact order it was originally
#    authored, so new shared strings land at the expected indices.
# ---------------------------------------------------------------------

$ws.Range("D52").Value2 = "Create, Edit for Auth"
$ws.Range("D53").Value2 = "Login Method"
$ws.Range("D54").Value2 = "Correcting every method and form"
$ws.Range("D56").Value2 = "Adapting every request for token auth"
$ws.Range("D58").Value2 = "GenerateController"
$ws.Range("D59").Value2 = "Generate Component"
$ws.Range("D55").Value2 = "Router - AuthGuards"
$ws.Range("D57").Value2 = "Backend controllers"
$ws.Range("D60").Value2 = "About page and documentation"

# ---------------------------------------------------------------------
# 3. Update the view: drop the frozen topLeftCell and move the selection
# ---------------------------------------------------------------------

$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D60").Select() | Out-Null
